# "update gebied onbekend en uitrustingsgraad/niveau"
#
# Insert a new row 2 ("gebied onbekend") above the existing data, pushing
# the current rows 2-27 down to rows 3-28. The new row records a #NULL!
# error code in column A (arrondiss2018) and 99993 in column B (provincie).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows down by inserting a blank row at row 2.
$ws.Rows("2:2").Insert()

# B2 should look like the rest of column B (same number format as its
# neighbours) - grab that formatting from the cell right below it.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# New values for the inserted row.
$ws.Range("A2").Value = "#NULL!"
$ws.Range("B2").Value = 99993

# A2 carries the error code, formatted as a percentage (distinct from the
# plain-integer style used by the rest of column A).
$ws.Range("A2").NumberFormat = "0.00%"
